$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 2945
$ws.Range("L3").Value = 2976
$ws.Range("L4").Value = 785
$ws.Range("L5").Value = 169
$ws.Range("L6").Value = 2670
$ws.Range("L7").Value = 9545

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L3").Value = 9
$ws.Range("L4").Value = 36
$ws.Range("L6").Value = 75
$ws.Range("L7").Value = 323
$ws.Range("L8").Value = 607
$ws.Range("L11").Value = 163
$ws.Range("L13").Value = 14
$ws.Range("L15").Value = 70
$ws.Range("L18").Value = 67
$ws.Range("L19").Value = 270
$ws.Range("L20").Value = 243
$ws.Range("L23").Value = 98
$ws.Range("L27").Value = 92
$ws.Range("L29").Value = 508
$ws.Range("L33").Value = 438
$ws.Range("L34").Value = 59
$ws.Range("L36").Value = 132
$ws.Range("L37").Value = 353
$ws.Range("L41").Value = 44
$ws.Range("L43").Value = 75
$ws.Range("L47").Value = 76
$ws.Range("L48").Value = 129
$ws.Range("L52").Value = 193
$ws.Range("L55").Value = 89
$ws.Range("L56").Value = 4
$ws.Range("L59").Value = 14
$ws.Range("L63").Value = 30
$ws.Range("L64").Value = 61
$ws.Range("L65").Value = 174
$ws.Range("L66").Value = 23
$ws.Range("L67").Value = 352
$ws.Range("L72").Value = 47
$ws.Range("L73").Value = 82
$ws.Range("L75").Value = 36
$ws.Range("L76").Value = 124
$ws.Range("L78").Value = 118
$ws.Range("L83").Value = 230
$ws.Range("L87").Value = 31
$ws.Range("L89").Value = 124
$ws.Range("L95").Value = 127
$ws.Range("L97").Value = 85
$ws.Range("L99").Value = 158
$ws.Range("L101").Value = 9545

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 98
$ws.Range("L3").Value = 100
$ws.Range("L7").Value = 323

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L4").Value = 13
$ws.Range("L7").Value = 163

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L2").Value = 40
$ws.Range("L3").Value = 35
$ws.Range("L7").Value = 124

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L3").Value = 55
$ws.Range("L7").Value = 193

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 173
$ws.Range("L3").Value = 204
$ws.Range("L6").Value = 166
$ws.Range("L7").Value = 607

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 76
$ws.Range("L3").Value = 91
$ws.Range("L6").Value = 50
$ws.Range("L7").Value = 230

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 124
$ws.Range("L6").Value = 152
$ws.Range("L7").Value = 438

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L3").Value = 38
$ws.Range("L7").Value = 127

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 104
$ws.Range("L3").Value = 103
$ws.Range("L5").Value = 13
$ws.Range("L7").Value = 353

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 63
$ws.Range("L6").Value = 46
$ws.Range("L7").Value = 174

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L6").Value = 40
$ws.Range("L7").Value = 158

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 107
$ws.Range("L7").Value = 352

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 157
$ws.Range("L3").Value = 189
$ws.Range("L6").Value = 134
$ws.Range("L7").Value = 508

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L3").Value = 28
$ws.Range("L7").Value = 129

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L3").Value = 82
$ws.Range("L7").Value = 270

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L3").Value = 19
$ws.Range("L7").Value = 124

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L3").Value = 22
$ws.Range("L7").Value = 75

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("L3").Value = 15
$ws.Range("L7").Value = 44

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("L5").Value = 6
$ws.Range("L6").Value = 14

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L3").Value = 34
$ws.Range("L7").Value = 118

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L2").Value = 33
$ws.Range("L7").Value = 89

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 98

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L3").Value = 14
$ws.Range("L4").Value = 6
$ws.Range("L7").Value = 61

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 77
$ws.Range("L7").Value = 243

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L2").Value = 25
$ws.Range("L7").Value = 67

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L2").Value = 52
$ws.Range("L7").Value = 132

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L2").Value = 15
$ws.Range("L7").Value = 59

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L3").Value = 27
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 76

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 27
$ws.Range("L7").Value = 70

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 23

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L6").Value = 20
$ws.Range("L7").Value = 82

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("L2").Value = 6
$ws.Range("L7").Value = 14

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L6").Value = 46
$ws.Range("L7").Value = 85

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L4").Value = 13
$ws.Range("L7").Value = 92

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("L2").Value = 18
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L3").Value = 22
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 75

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("L3").Value = 10
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 47

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("L2").Value = 1
$ws.Range("L7").Value = 4

$ws = $wb.Worksheets.Item("Andersonville")
$ws.Range("L2").Value = 4
$ws.Range("L7").Value = 9

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("L3").Value = 11
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("L3").Value = 4
$ws.Range("L7").Value = 31
